$wb = $excel.ActiveWorkbook

$wsVar = $wb.Worksheets.Item("Variables")
$wsCat = $wb.Worksheets.Item("Categories")

# Remove the "ff10" / "Total number of stillbirths" row (row 20) from the Variables sheet.
# This shifts all subsequent rows up by one.
$wsVar.Rows.Item(20).Delete()

# Update selections / active sheet state to match the target file.
$wsCat.Activate()
$wsCat.Range("I48").Select()
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$wsVar.Activate()
$wsVar.Range("E20").Select()

$wb.Save()
